$wb = $excel.ActiveWorkbook

# --- Add a new "Speaker" worksheet, positioned right before "Cookies" ---
$cookies = $wb.Worksheets.Item("Cookies")
$speaker = $wb.Worksheets.Add($cookies)
$speaker.Name = "Speaker"

# Column widths (approximate the bestFit widths used on the similar TV sheet)
$speaker.Columns.Item(5).ColumnWidth = 19.8320625
$speaker.Columns.Item(6).ColumnWidth = 15.168
$speaker.Columns.Item(7).ColumnWidth = 10.8320625

# Header row
$speaker.Range("E10").Value = "Device"
$speaker.Range("F10").Value = "Installed base (M)"
$speaker.Range("G10").Value = "Power Draw (w)"
$speaker.Range("H10").Value = "Usage (h/day)"

# Smart speaker device defaults
$speaker.Range("E11").Value = "Google Home Mini"
$speaker.Range("F11").Value = 4
$speaker.Range("G11").Value = 1.7
$speaker.Range("H11").Value = 3.5

$speaker.Range("E12").Value = "Amazon Echo (2nd gen)"
$speaker.Range("F12").Value = 35
$speaker.Range("G12").Value = 2.4
$speaker.Range("H12").Value = 3.5

$speaker.Range("E13").Value = "Google Home"
$speaker.Range("F13").Value = 8
$speaker.Range("G13").Value = 2.2000000000000002
$speaker.Range("H13").Value = 3.5

$speaker.Range("E14").Value = "Apple HomePod"
$speaker.Range("F14").Value = 3
$speaker.Range("G14").Value = 5.9
$speaker.Range("H14").Value = 3.5

# Blank helper cell keeps the same "0.0" number format as the summary cell below
$speaker.Range("G16").NumberFormat = "0.0"

# Weighted-average power draw summary
$speaker.Range("E18").Value = "Smart Speaker"
$speaker.Range("G18").NumberFormat = "0.0"
$speaker.Range("G18").Formula = "=SUMPRODUCT(F11:F16,G11:G16)/SUM(F11:F14)"

# Match the saved selection/active-cell state of the new sheet
$speaker.Range("H18").Select()
